# Update cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.863.23"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.742.91"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.08"
$ws.Range("E5").Value = "  -5.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5150"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2805"
$ws.Range("E8").Value = "  +6.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.13"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06088"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "1.747.28"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06959"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.21"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6307"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.483"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.31"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "25.882.07"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006574"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").Value = "1.965.95"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.084"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.394"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.109"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.27"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.505"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.818"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.98"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.40"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08274"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.604"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.393"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04369"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.621"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6025"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.667"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01548"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.896"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.28"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3811"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7179"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.895"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.253"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1099"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.20"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.61"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.07%  "
